$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update dimension via header + value edits (Excel recalculates dimension automatically based on content)

# 2. Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3. Title-case the Spanish connector words (de/del/la/el/los/las/y) in state/municipality names
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B7").Value = "San Francisco De Los Romo"
$ws.Range("B21").Value = "Amatenango De La Frontera"
$ws.Range("B22").Value = "Amatenango Del Valle"
$ws.Range("B25").Value = "Bejucal De Ocampo"
$ws.Range("B30").Value = "Chiapa De Corzo"
$ws.Range("B35").Value = "Comitán De Domínguez"
$ws.Range("B54").Value = "Marqués De Comillas"
$ws.Range("B55").Value = "Mazapa De Madero"
$ws.Range("B60").Value = "Ocozocoautla De Espinosa"
$ws.Range("B66").Value = "San Cristóbal De Las Casas"
$ws.Range("B90").Value = "Hidalgo Del Parral"
$ws.Range("B102").Value = "San Juan De Sabinas"
$ws.Range("A113").Value = "Ciudad De México"
$ws.Range("B135").Value = "Nombre De Dios"
$ws.Range("B137").Value = "San Juan De Guadalupe"
$ws.Range("A143").Value = "Estado De México"
$ws.Range("B143").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B146").Value = "Almoloya De Alquisiras"
$ws.Range("B147").Value = "Almoloya Del Río"
$ws.Range("B149").Value = "Atizapán De Zaragoza"
$ws.Range("B158").Value = "Ecatepec De Morelos"
$ws.Range("B161").Value = "Ixtapan De La Sal"
$ws.Range("B165").Value = "Naucalpan De Juárez"
$ws.Range("B170").Value = "San Felipe Del Progreso"
$ws.Range("B171").Value = "San Martín De Las Pirámides"
$ws.Range("B177").Value = "Tenango Del Valle"
$ws.Range("B183").Value = "Tlalnepantla De Baz"
$ws.Range("B187").Value = "Valle De Bravo"
$ws.Range("B188").Value = "Valle De Chalco Solidaridad"
$ws.Range("B196").Value = "Apaseo El Alto"
$ws.Range("B202").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B214").Value = "San Diego De La Unión"
$ws.Range("B216").Value = "San Francisco Del Rincón"
$ws.Range("B218").Value = "San Luis De La Paz"
$ws.Range("B219").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B223").Value = "Valle De Santiago"
$ws.Range("B228").Value = "Acapulco De Juárez"
$ws.Range("B232").Value = "Atlamajalcingo Del Monte"
$ws.Range("B234").Value = "Atoyac De Álvarez"
$ws.Range("B235").Value = "Ayutla De Los Libres"
$ws.Range("B237").Value = "Chilapa De Álvarez"
$ws.Range("B238").Value = "Chilpancingo De Los Bravo"
$ws.Range("B241").Value = "Coyuca De Benítez"
$ws.Range("B242").Value = "Coyuca De Catalán"
$ws.Range("B249").Value = "Huitzuco De Los Figueroa"
$ws.Range("B250").Value = "Iguala De La Independencia"
$ws.Range("B251").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B252").Value = "Zihuatanejo De Azueta"
$ws.Range("B266").Value = "Taxco De Alarcón"
$ws.Range("B268").Value = "Técpan De Galeana"
$ws.Range("B270").Value = "Tepecoacuilco De Trujano"
$ws.Range("B278").Value = "Agua Blanca De Iturbide"
$ws.Range("B279").Value = "Atotonilco El Grande"
$ws.Range("B284").Value = "Cuautepec De Hinojosa"
$ws.Range("B288").Value = "Huejutla De Reyes"
$ws.Range("B291").Value = "Jacala De Ledezma"
$ws.Range("B297").Value = "Mixquiahuala De Juárez"
$ws.Range("B298").Value = "Pachuca De Soto"
$ws.Range("B299").Value = "Progreso De Obregón"
$ws.Range("B302").Value = "Santiago De Anaya"
$ws.Range("B305").Value = "Tenango De Doria"
$ws.Range("B307").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B309").Value = "Tezontepec De Aldama"
$ws.Range("B314").Value = "Tula De Allende"
$ws.Range("B315").Value = "Tulancingo De Bravo"
$ws.Range("B322").Value = "Autlán De Navarro"
$ws.Range("B327").Value = "Huejuquilla El Alto"
$ws.Range("B331").Value = "Lagos De Moreno"
$ws.Range("B334").Value = "Ojuelos De Jalisco"
$ws.Range("B338").Value = "San Juan De Los Lagos"
$ws.Range("B339").Value = "San Miguel El Alto"
$ws.Range("B340").Value = "San Sebastián Del Oeste"
$ws.Range("B343").Value = "Tlajomulco De Zúñiga"
$ws.Range("B346").Value = "Unión De San Antonio"
$ws.Range("B363").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B364").Value = "Cojumatlán De Régules"
$ws.Range("B422").Value = "Puente De Ixtla"
$ws.Range("B424").Value = "Tlaltizapán De Zapata"
$ws.Range("B432").Value = "Ixtlán Del Río"
$ws.Range("B448").Value = "Mier Y Noriega"
$ws.Range("B450").Value = "San Nicolás De Los Garza"
$ws.Range("B452").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B457").Value = "Constancia Del Rosario"
$ws.Range("B460").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B461").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B462").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B463").Value = "Ixtlán De Juárez"
$ws.Range("B464").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B469").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B470").Value = "Oaxaca De Juárez"
$ws.Range("B471").Value = "Ocotlán De Morelos"
$ws.Range("B472").Value = "Pinotepa De Don Luis"
$ws.Range("B474").Value = "Putla Villa De Guerrero"
$ws.Range("B504").Value = "San Miguel Del Puerto"
$ws.Range("B556").Value = "Santo Domingo De Morelos"
$ws.Range("B562").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B563").Value = "Tataltepec De Valdés"
$ws.Range("B564").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B565").Value = "Tlacolula De Matamoros"
$ws.Range("B566").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B567").Value = "Villa Sola De Vega"
$ws.Range("B568").Value = "Zapotitlán Del Río"
$ws.Range("B588").Value = "Huehuetlán El Chico"
$ws.Range("B589").Value = "Ixcamilpa De Guerrero"
$ws.Range("B590").Value = "Izúcar De Matamoros"
$ws.Range("B592").Value = "Los Reyes De Juárez"
$ws.Range("B594").Value = "Palmar De Bravo"
$ws.Range("B602").Value = "San Salvador El Seco"
$ws.Range("B608").Value = "Tepexi De Rodríguez"
$ws.Range("B609").Value = "Tetela De Ocampo"
$ws.Range("B611").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B623").Value = "Amealco De Bonfil"
$ws.Range("B625").Value = "Cadereyta De Montes"
$ws.Range("B627").Value = "Jalpan De Serra"
$ws.Range("B629").Value = "San Juan Del Río"
$ws.Range("B641").Value = "Ciudad Del Maíz"
$ws.Range("B654").Value = "Santa María Del Río"
$ws.Range("B656").Value = "Soledad De Graciano Sánchez"
$ws.Range("B662").Value = "Villa De Arista"
$ws.Range("B663").Value = "Villa De Guadalupe"
$ws.Range("B664").Value = "Villa De La Paz"
$ws.Range("B665").Value = "Villa De Ramos"
$ws.Range("B689").Value = "Jalpa De Méndez"
$ws.Range("B707").Value = "Soto La Marina"
$ws.Range("B714").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B721").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B730").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B734").Value = "Amatlán De Los Reyes"
$ws.Range("B739").Value = "Boca Del Río"
$ws.Range("B750").Value = "Cosamaloapan De Carpio"
$ws.Range("B762").Value = "Hueyapan De Ocampo"
$ws.Range("B763").Value = "Ignacio De La Llave"
$ws.Range("B765").Value = "Ixhuatlán Del Café"
$ws.Range("B771").Value = "Juchique De Ferrer"
$ws.Range("B776").Value = "Martínez De La Torre"
$ws.Range("B785").Value = "Ozuluama De Mascareñas"
$ws.Range("B788").Value = "Paso Del Macho"
$ws.Range("B790").Value = "Poza Rica De Hidalgo"
$ws.Range("B796").Value = "Sayula De Alemán"
$ws.Range("B797").Value = "Soledad De Doblado"
$ws.Range("B814").Value = "Tlacotepec De Mejía"
$ws.Range("B821").Value = "Vega De Alatorre"
$ws.Range("B835").Value = "Cañitas De Felipe Pescador"
$ws.Range("B845").Value = "Nochistlán De Mejía"
$ws.Range("B846").Value = "Noria De Ángeles"
$ws.Range("B855").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B858").Value = "Villa De Cos"

# 4. Fix tiny floating point rounding for D71
$ws.Range("D71").Value = 0.009575104727707961

# 5. Remove trailing footnote rows 863-867 (they fall outside the data range now)
$ws.Range("A863:D867").EntireRow.Delete() | Out-Null

